$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 42321
$ws.Range("E2").Value = 830492108062
$ws.Range("F2").Value = 13695151538
$ws.Range("G2").Value = 1.30507

$ws.Range("D3").Value = 2271.71
$ws.Range("E3").Value = 273218759616
$ws.Range("F3").Value = 6117519730
$ws.Range("G3").Value = 0.17584

$ws.Range("D4").Value = 0.999547
$ws.Range("E4").Value = 95993219168
$ws.Range("F4").Value = 22257190702
$ws.Range("G4").Value = 0.0375

$ws.Range("D5").Value = 305.62
$ws.Range("E5").Value = 47075580775
$ws.Range("F5").Value = 481764802
$ws.Range("G5").Value = 0.44138

$ws.Range("D6").Value = 97.17
$ws.Range("E6").Value = 42173109943
$ws.Range("F6").Value = 1951515852
$ws.Range("G6").Value = 5.75738

$ws.Range("D7").Value = 0.5288
$ws.Range("E7").Value = 28734960259
$ws.Range("F7").Value = 419890891
$ws.Range("G7").Value = -0.32527

$ws.Range("D8").Value = 0.999524
$ws.Range("E8").Value = 26108619516
$ws.Range("F8").Value = 3170243521
$ws.Range("G8").Value = -0.08237999999999999

$ws.Range("D9").Value = 2272.23
$ws.Range("E9").Value = 21323826336
$ws.Range("F9").Value = 15608200
$ws.Range("G9").Value = 0.34966

$ws.Range("D10").Value = 0.493538
$ws.Range("E10").Value = 17373390598
$ws.Range("F10").Value = 302844795
$ws.Range("G10").Value = 2.07222

$ws.Range("B11").Value = "AVAX"
$ws.Range("C11").Value = "Avalanche"
$ws.Range("D11").Value = 35.78
$ws.Range("E11").Value = 13178504751
$ws.Range("F11").Value = 779976454
$ws.Range("G11").Value = 10.41223

$ws.Range("B12").Value = "DOGE"
$ws.Range("C12").Value = "Dogecoin"
$ws.Range("D12").Value = 0.079541
$ws.Range("E12").Value = 11395278106
$ws.Range("F12").Value = 281364102
$ws.Range("G12").Value = -0.37941

$ws.Range("D13").Value = 0.111727
$ws.Range("E13").Value = 9846267907
$ws.Range("F13").Value = 245437131
$ws.Range("G13").Value = -2.41961

$ws.Range("D14").Value = 6.67
$ws.Range("E14").Value = 8833009018
$ws.Range("F14").Value = 160313695
$ws.Range("G14").Value = 0.22773

$ws.Range("D15").Value = 14.45
$ws.Range("E15").Value = 8222508628
$ws.Range("F15").Value = 297341015
$ws.Range("G15").Value = 1.49187

$ws.Range("B16").Value = "MATIC"
$ws.Range("C16").Value = "Polygon"
$ws.Range("D16").Value = 0.792156
$ws.Range("E16").Value = 7362214771
$ws.Range("F16").Value = 310379960
$ws.Range("G16").Value = 3.07889

$ws.Range("B17").Value = "TON"
$ws.Range("C17").Value = "Toncoin"
$ws.Range("D17").Value = 2.11
$ws.Range("E17").Value = 7319162621
$ws.Range("F17").Value = 12450200
$ws.Range("G17").Value = 2.15062

$ws.Range("D18").Value = 42165
$ws.Range("E18").Value = 6659446151
$ws.Range("F18").Value = 155470938
$ws.Range("G18").Value = 1.13858

$ws.Range("B19").Value = "ICP"
$ws.Range("C19").Value = "Internet Computer"
$ws.Range("D19").Value = 12.49
$ws.Range("E19").Value = 5722725459
$ws.Range("F19").Value = 139861241
$ws.Range("G19").Value = -0.55908

$ws.Range("B20").Value = "SHIB"
$ws.Range("C20").Value = "Shiba Inu"
$ws.Range("D20").Value = 0.000009099999999999999
$ws.Range("E20").Value = 5364989640
$ws.Range("F20").Value = 80224391
$ws.Range("G20").Value = 0.55628

$ws.Range("D21").Value = 0.998708
$ws.Range("E21").Value = 5227568119
$ws.Range("F21").Value = 140419088
$ws.Range("G21").Value = -0.15169

$ws.Range("B22").Value = "LTC"
$ws.Range("C22").Value = "Litecoin"
$ws.Range("D22").Value = 67.45999999999999
$ws.Range("E22").Value = 5004501090
$ws.Range("F22").Value = 209212851
$ws.Range("G22").Value = 0.5116000000000001

$ws.Range("B23").Value = "BCH"
$ws.Range("C23").Value = "Bitcoin Cash"
$ws.Range("D23").Value = 240.82
$ws.Range("E23").Value = 4737359352
$ws.Range("F23").Value = 135652266
$ws.Range("G23").Value = 0.18007

$ws.Range("B24").Value = "UNI"
$ws.Range("C24").Value = "Uniswap"
$ws.Range("D24").Value = 5.97
$ws.Range("E24").Value = 4509306868
$ws.Range("F24").Value = 61099003
$ws.Range("G24").Value = 0.66314

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "LEO Token"
$ws.Range("D25").Value = 4.02
$ws.Range("E25").Value = 3730353999
$ws.Range("F25").Value = 979564
$ws.Range("G25").Value = -1.3253

$ws.Range("B26").Value = "ATOM"
$ws.Range("C26").Value = "Cosmos Hub"
$ws.Range("D26").Value = 9.52
$ws.Range("E26").Value = 3660579120
$ws.Range("F26").Value = 139692498
$ws.Range("G26").Value = -0.06261

$ws.Range("D27").Value = 23.75
$ws.Range("E27").Value = 3410126213
$ws.Range("F27").Value = 143586131
$ws.Range("G27").Value = -0.86826

$ws.Range("B28").Value = "INJ"
$ws.Range("C28").Value = "Injective"
$ws.Range("D28").Value = 37.58
$ws.Range("E28").Value = 3329054311
$ws.Range("F28").Value = 188625946
$ws.Range("G28").Value = 7.6221

$ws.Range("B29").Value = "XLM"
$ws.Range("C29").Value = "Stellar"
$ws.Range("D29").Value = 0.11462
$ws.Range("E29").Value = 3255054834
$ws.Range("F29").Value = 50189685
$ws.Range("G29").Value = -1.6747

$ws.Range("B30").Value = "OKB"
$ws.Range("C30").Value = "OKB"
$ws.Range("D30").Value = 52.04
$ws.Range("E30").Value = 3128171189
$ws.Range("F30").Value = 8881733
$ws.Range("G30").Value = -2.38593

$ws.Range("B31").Value = "NEAR"
$ws.Range("C31").Value = "NEAR Protocol"
$ws.Range("D31").Value = 2.97
$ws.Range("E31").Value = 3051305138
$ws.Range("F31").Value = 197751855
$ws.Range("G31").Value = 2.8764

$ws.Range("B32").Value = "APT"
$ws.Range("C32").Value = "Aptos"
$ws.Range("D32").Value = 8.949999999999999
$ws.Range("E32").Value = 3025189835
$ws.Range("F32").Value = 109030465
$ws.Range("G32").Value = -0.94171

$ws.Range("B33").Value = "OP"
$ws.Range("C33").Value = "Optimism"
$ws.Range("D33").Value = 3.09
$ws.Range("E33").Value = 2959472583
$ws.Range("F33").Value = 138423608
$ws.Range("G33").Value = 1.35537

$ws.Range("B34").Value = "XMR"
$ws.Range("C34").Value = "Monero"
$ws.Range("D34").Value = 159.89
$ws.Range("E34").Value = 2904573106
$ws.Range("F34").Value = 50570719
$ws.Range("G34").Value = -0.68785

$ws.Range("B35").Value = "LDO"
$ws.Range("C35").Value = "Lido DAO"
$ws.Range("D35").Value = 3.11
$ws.Range("E35").Value = 2777486006
$ws.Range("F35").Value = 47030840
$ws.Range("G35").Value = 3.14679

$ws.Range("B36").Value = "TIA"
$ws.Range("C36").Value = "Celestia"
$ws.Range("D36").Value = 17.04
$ws.Range("E36").Value = 2741641285
$ws.Range("F36").Value = 118633054
$ws.Range("G36").Value = 1.26898

$ws.Range("B37").Value = "IMX"
$ws.Range("C37").Value = "Immutable"
$ws.Range("D37").Value = 1.94
$ws.Range("E37").Value = 2624209123
$ws.Range("F37").Value = 49551442
$ws.Range("G37").Value = 0.77366

$ws.Range("B38").Value = "FIL"
$ws.Range("C38").Value = "Filecoin"
$ws.Range("D38").Value = 5.21
$ws.Range("E38").Value = 2612613145
$ws.Range("F38").Value = 121372607
$ws.Range("G38").Value = -1.3966

$ws.Range("B39").Value = "FDUSD"
$ws.Range("C39").Value = "First Digital USD"
$ws.Range("D39").Value = 0.998757
$ws.Range("E39").Value = 2588859435
$ws.Range("F39").Value = 2065762098
$ws.Range("G39").Value = -0.03568

$ws.Range("B40").Value = "HBAR"
$ws.Range("C40").Value = "Hedera"
$ws.Range("D40").Value = 0.07421999999999999
$ws.Range("E40").Value = 2502058238
$ws.Range("F40").Value = 23179653
$ws.Range("G40").Value = -0.1024

$ws.Range("B41").Value = "KAS"
$ws.Range("C41").Value = "Kaspa"
$ws.Range("D41").Value = 0.105592
$ws.Range("E41").Value = 2393681407
$ws.Range("F41").Value = 14871555
$ws.Range("G41").Value = 0.89373

$ws.Range("B42").Value = "ARB"
$ws.Range("C42").Value = "Arbitrum"
$ws.Range("D42").Value = 1.84
$ws.Range("E42").Value = 2350679739
$ws.Range("F42").Value = 356313796
$ws.Range("G42").Value = 2.5058

$ws.Range("B43").Value = "TAO"
$ws.Range("C43").Value = "Bittensor"
$ws.Range("D43").Value = 372.09
$ws.Range("E43").Value = 2271967245
$ws.Range("F43").Value = 11501229
$ws.Range("G43").Value = 11.69347

$ws.Range("B44").Value = "STX"
$ws.Range("C44").Value = "Stacks"
$ws.Range("D44").Value = 1.52
$ws.Range("E44").Value = 2191258743
$ws.Range("F44").Value = 44074928
$ws.Range("G44").Value = 0.21957

$ws.Range("B45").Value = "CRO"
$ws.Range("C45").Value = "Cronos"
$ws.Range("D45").Value = 0.081175
$ws.Range("E45").Value = 2154171967
$ws.Range("F45").Value = 5721873
$ws.Range("G45").Value = -0.37099

$ws.Range("D46").Value = 0.02842794
$ws.Range("E46").Value = 2069582345
$ws.Range("F46").Value = 36288019
$ws.Range("G46").Value = 2.05685

$ws.Range("B47").Value = "MNT"
$ws.Range("C47").Value = "Mantle"
$ws.Range("D47").Value = 0.648046
$ws.Range("E47").Value = 2054400826
$ws.Range("F47").Value = 45880692
$ws.Range("G47").Value = -0.9555399999999999

$ws.Range("D48").Value = 1988.7
$ws.Range("E48").Value = 1837061478
$ws.Range("F48").Value = 56147238
$ws.Range("G48").Value = -1.62676

$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "Sei"
$ws.Range("D49").Value = 0.655002
$ws.Range("E49").Value = 1615929947
$ws.Range("F49").Value = 157221549
$ws.Range("G49").Value = 2.23989

$ws.Range("B50").Value = "TUSD"
$ws.Range("C50").Value = "TrueUSD"
$ws.Range("D50").Value = 0.985966
$ws.Range("E50").Value = 1595947812
$ws.Range("F50").Value = 54959830
$ws.Range("G50").Value = -0.22305

$ws.Range("B51").Value = "QNT"
$ws.Range("C51").Value = "Quant"
$ws.Range("D51").Value = 107.82
$ws.Range("E51").Value = 1569694704
$ws.Range("F51").Value = 17163814
$ws.Range("G51").Value = -1.30366

